# Actualizacion Certificacion, arreglo de descarga de detalle
#
# 1. B3's label changes from "PAGO TELGUA EN LINEA" to
#    "CxCAjena - PAGO TELGUA EN LINEA" (the old shared string becomes
#    unused and drops out of the table; the new text is appended).
# 2. The A3:A4 merge (the helper "Rango" layout) is removed.
# 3. The now-unmerged A3:A4 cells keep vertical centering but lose the
#    horizontal centering that the merged layout had.
# 4. The active selection ends up on B4 instead of A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the label text.
$ws.Range("B3").Value = "CxCAjena - PAGO TELGUA EN LINEA"

# 2) Remove the A3:A4 merge.
[void]$ws.Range("A3:A4").UnMerge()

# 3) Drop horizontal centering on A3:A4 (vertical centering is kept as-is).
$ws.Range("A3:A4").HorizontalAlignment = 1

# 4) Move the selection to B4.
[void]$ws.Range("B4").Select()
